$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71, column A currently holds the phone number 71717173 stored as
# text (it should have been numeric, like every other row in the sheet).
# Fix the type by writing a genuine numeric literal.
$ws.Cells.Item(71, 1).Value = 71717173

# Append the new payment record as row 72 (phone 76442711, Cash, $408,
# no discount/birthday-discount/points/reward applied).
#
# Phone numbers in this sheet are stored as text, and a couple of the
# numeric columns are left as blank text cells rather than true numbers
# when there's no discount data. Prefixing a value with an apostrophe
# forces Excel to store it as text instead of auto-coercing it to a
# Number (this is also how the blank-but-present text cells are
# produced: an apostrophe with nothing after it). Resetting the style to
# Normal afterwards strips the implicit "Text" number format Excel
# applies when it protects a value from numeric coercion, so the cell
# ends up with the same default styling as its neighbours.
$ws.Cells.Item(72, 1).Value = "'76442711"
$ws.Cells.Item(72, 1).Style = "Normal"

$ws.Cells.Item(72, 2).Value = "'"
$ws.Cells.Item(72, 2).Style = "Normal"

$ws.Cells.Item(72, 3).Value = "Cash"
$ws.Cells.Item(72, 4).Value = "2025-08-20T08:42:01"
$ws.Cells.Item(72, 5).Value = 408

$ws.Cells.Item(72, 6).Value = "'"
$ws.Cells.Item(72, 6).Style = "Normal"

$ws.Cells.Item(72, 7).Value = 408
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 0
